$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 109853088
$ws.Range("B9").Value = 73698
$ws.Range("E9").Value = 1467
$ws.Range("F9").Value = 'Rödbrun blekspik'
$ws.Range("G9").Value = 'Sclerophora coniophaea'
$ws.Range("H9").Value = '(Norman) J.Mattsson & Middelb.'
$ws.Range("Q9").Value = 458682.9831145834
$ws.Range("R9").Value = 6910834.730994376
$ws.Range("A10").Value = 109852962
$ws.Range("B10").Value = 89406
$ws.Range("E10").Value = 1204
$ws.Range("F10").Value = 'Gränsticka'
$ws.Range("G10").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H10").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("Q10").Value = 458604.1467445081
$ws.Range("R10").Value = 6910745.989096207
$ws.Range("A11").Value = 109852942
$ws.Range("B11").Value = 89673
$ws.Range("E11").Value = 658
$ws.Range("F11").Value = 'Rosenticka'
$ws.Range("G11").Value = 'Rhodofomes roseus'
$ws.Range("H11").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q11").Value = 458624.8446887356
$ws.Range("R11").Value = 6910763.397753252
$ws.Range("A12").Value = 109876605
$ws.Range("B12").Value = 56395
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("M12").Value = 'äldre spår'
$ws.Range("Q12").Value = 458587.280630524
$ws.Range("R12").Value = 6910812.67010971
$ws.Range("S12").Value = 48
$ws.Range("AW12").Value = 'Martin Kämpedal'
$ws.Range("AX12").Value = 'Martin Kämpedal'
$ws.Range("A13").Value = 109875646
$ws.Range("B13").Value = 77507
$ws.Range("E13").Value = 230405
$ws.Range("F13").Value = 'Garnlav (ssp. sarmentosa)'
$ws.Range("G13").Value = 'Alectoria sarmentosa subsp. sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("J13").Value = ""
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = ""
$ws.Range("P13").Value = 'Torpvallen, Hjd'
$ws.Range("Q13").Value = 458548.9991899654
$ws.Range("R13").Value = 6910766.194164375
$ws.Range("S13").Value = 10
$ws.Range("AF13").Value = ""
$ws.Range("A14").Value = 109875640
$ws.Range("Q14").Value = 458575.8866563534
$ws.Range("R14").Value = 6910681.723629645
$ws.Range("A15").Value = 109875687
$ws.Range("Q15").Value = 459008.9636710359
$ws.Range("R15").Value = 6910920.905611085
$ws.Range("A16").Value = 109875690
$ws.Range("Q16").Value = 458931.8478696992
$ws.Range("R16").Value = 6911009.699956965
$ws.Range("A17").Value = 109875659
$ws.Range("P17").Value = 'Torpvallen, Hjd'
$ws.Range("Q17").Value = 458700.9619114345
$ws.Range("R17").Value = 6910858.215982608
$ws.Range("A18").Value = 109875679
$ws.Range("Q18").Value = 458779.0883940465
$ws.Range("R18").Value = 6911003.208538387
$ws.Range("A19").Value = 109875671
$ws.Range("Q19").Value = 458558.4460608965
$ws.Range("R19").Value = 6910890.18893374
$ws.Range("A20").Value = 109875649
$ws.Range("Q20").Value = 458670.7282092119
$ws.Range("R20").Value = 6910784.679410813
$ws.Range("A21").Value = 109875664
$ws.Range("Q21").Value = 458639.1266098225
$ws.Range("R21").Value = 6910902.207586206
$ws.Range("A22").Value = 109875775
$ws.Range("Q22").Value = 458881.5673188641
$ws.Range("R22").Value = 6910933.156907165
$ws.Range("A23").Value = 109875684
$ws.Range("Q23").Value = 459003.819864114
$ws.Range("R23").Value = 6910842.877383335
$ws.Range("A24").Value = 109853168
$ws.Range("B24").Value = 77588
$ws.Range("E24").Value = 864
$ws.Range("F24").Value = 'Knottrig blåslav'
$ws.Range("G24").Value = 'Hypogymnia bitteri'
$ws.Range("H24").Value = '(Lynge) Ahti'
$ws.Range("J24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("P24").Value = 'Skedflötarna, Hjd'
$ws.Range("Q24").Value = 458567.6177838827
$ws.Range("R24").Value = 6910615.35053566
$ws.Range("S24").Value = 25
$ws.Range("AF24").ClearContents()
$ws.Range("AW24").Value = 'lennart karlsson'
$ws.Range("AX24").Value = 'lennart karlsson'
$ws.Range("A25").Value = 109875605
$ws.Range("P25").Value = 'Garnlav, Hjd'
$ws.Range("Q25").Value = 458561.882745445
$ws.Range("R25").Value = 6910603.334865634
$ws.Range("A26").Value = 109875589
$ws.Range("Q26").Value = 458504.9961458603
$ws.Range("R26").Value = 6910520.826033623
